$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F18").Value = 0.022958220000000001
$ws.Range("F19").Value = 0.3307466
$ws.Range("F20").Value = 0.069393700000000003
$ws.Range("F21").Value = 0.57690149999999996
$ws.Range("F22").Value = 0.024115020000000001
$ws.Range("F23").Value = 0.34303359999999999
$ws.Range("F24").Value = 0.071385560000000001
$ws.Range("F25").Value = 0.56146580000000001
$ws.Range("F26").Value = 0.013475050000000001
$ws.Range("F27").Value = 0.21222250000000001
$ws.Range("F28").Value = 0.04450759
$ws.Range("F29").Value = 0.72979490000000002
$ws.Range("F30").Value = 0.01495052
$ws.Range("F31").Value = 0.21552859999999999
$ws.Range("F32").Value = 0.043749719999999999
$ws.Range("F33").Value = 0.72577119999999995

$ws.Range("F34").Select()
